# Auto-generated: update cryptocurrency price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.849.07'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -1.09%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.893.57'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -0.93%  '

# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7829'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -4.71%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '244.02'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.05%  '

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.05%  '

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -3.36%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '25.32'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -5.75%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07204'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.16%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08092'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.15%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.7650'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.99%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.479'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +3.47%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.897.32'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.91%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.37'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.79%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.157'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +4.13%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '29.843.37'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.10%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '13.99'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -1.75%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '243.66'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.94%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.000007794'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.07%  '

# Row 21
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.11%  '

# Row 22
$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.154.05'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.98%  '

# Row 23
$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = 'Chainlink'
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '8.147'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +15.35%  '

# Row 24
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.06%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1642'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.01%  '

# Row 26
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.02%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '163.04'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -2.53%  '

# Row 28
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -1.34%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.052'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -2.74%  '

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +2.85%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.550'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.42%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.503'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +4.55%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.124'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.42%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05567'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.93%  '

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -0.45%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7436'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +1.05%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.9978'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.32%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.617'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.16%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01920'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.27%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.782'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.52%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.147.78'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +14.09%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '73.97'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +0.66%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8514'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.18%  '

# Row 46
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.001'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.14%  '

# Row 47
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Quant'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '103.81'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.02%  '

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -1.65%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.970'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.52%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.468'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -1.65%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.995'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +9.82%  '
